$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows: lowercase names, swap/adjust balances
$ws.Range("A2").Value = "luiz"
$ws.Range("B2").Value = 1000

$ws.Range("A3").Value = "pedro"
$ws.Range("B3").Value = -1000

# Add new rows
$ws.Range("A4").Value = "joao"
$ws.Range("B4").Value = 1000

$ws.Range("A5").Value = "lucas"
$ws.Range("B5").Value = 1000
